# Update "want to go" counts (column F) on the "展览" and "全部类型" sheets
# to reflect a newer data snapshot, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 14486
$ws1.Range("F5").Value = 17173
$ws1.Range("F6").Value = 19
$ws1.Range("F8").Value = 49
$ws1.Range("F17").Value = 14
$ws1.Range("F20").Value = 1309
$ws1.Range("F25").Value = 7079
$ws1.Range("F27").Value = 32
$ws1.Range("F33").Value = 130
$ws1.Range("F36").Value = 4994

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 14486
$ws4.Range("F5").Value = 17174
$ws4.Range("F6").Value = 19
$ws4.Range("F8").Value = 49
$ws4.Range("F17").Value = 14
$ws4.Range("F20").Value = 1309
$ws4.Range("F26").Value = 7079
$ws4.Range("F28").Value = 32
$ws4.Range("F35").Value = 130
$ws4.Range("F38").Value = 4994
